$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DecisionTable")

# Row 28: new rule row appended to the decision table.
# A28 = "Test" (plain text)
$ws.Range("A28").Value = "Test"

# B28 = "60" stored as text (matches the existing numeric-looking text
# entries elsewhere in this CONDITION column, e.g. B20:B24).
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "60"
$ws.Range("B28").Style = "Normal"

# C28 = empty text value (same pattern as other blank CONDITION cells
# that still hold an empty shared string, e.g. A19/C20:C22).
$ws.Range("C28").Value = "'"
$ws.Range("C28").Style = "Normal"

# D28:F28 stay blank, but keep an explicit (styled) empty cell like the
# rest of the table.
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Style = "Normal"
$ws.Range("F28").Style = "Normal"

# G28 = "6" stored as text (ACTION column).
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "6"
$ws.Range("G28").Style = "Normal"
